# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio" (D) and "aragon" (G) columns of metadata get re-pointed
# to the newly curated sdmx/skos dimension vocabulary, and the obsolete
# "mapping-aragon.xlsx" reference (G5) is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (municipio-nombre block)
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column G (aragon block)
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("G4").Value = "URI-Comunidad"
$ws.Range("G5").ClearContents()
